$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Trey Murphy III", "SG,SF,PF", "New Orleans Pelicans"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Al Horford", "PF,C", "Boston Celtics"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
